$d = $word.ActiveDocument

$newParagraphsXml = @'
    <w:p>
      <w:r>
        <w:br w:type="page"/>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Segoe UI Emoji" w:hAnsi="Segoe UI Emoji" w:cs="Segoe UI Emoji"/>
        </w:rPr>
        <w:lastRenderedPageBreak/>
        <w:t>🏁🔥</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
        </w:rPr>
        <w:t>¡DIEGO FERNÁNDEZ GANA EN SILVERSTONE! UNA CARRERA DE F1 CON CHOQUES, DEBUTS Y PODIOS SORPRESA</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Segoe UI Emoji" w:hAnsi="Segoe UI Emoji" w:cs="Segoe UI Emoji"/>
        </w:rPr>
        <w:t>🔥🏁</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>Silverstone</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> fue escenario de una carrera de F1 con </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
        </w:rPr>
        <w:t>choques desde la primera curva, emociones en pista y debuts que dieron de qué hablar</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">. Todo arrancó con un homenaje al mítico Eddie </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>Jordan</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve">, donde el equipo </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>Jordan</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> y varios pilotos salieron con la skin de 1991.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Segoe UI Emoji" w:hAnsi="Segoe UI Emoji" w:cs="Segoe UI Emoji"/>
        </w:rPr>
        <w:t>🚀</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
        </w:rPr>
        <w:t>Diego Fernández salió desde la pole</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> y no soltó la punta: victoria sólida y sin errores.</w:t>
      </w:r>
      <w:r>
        <w:br/>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Segoe UI Emoji" w:hAnsi="Segoe UI Emoji" w:cs="Segoe UI Emoji"/>
        </w:rPr>
        <w:t>👶</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
        </w:rPr>
        <w:t xml:space="preserve">Hugo debutó con </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
        </w:rPr>
        <w:t>Kick</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
        </w:rPr>
        <w:t xml:space="preserve"> y se metió directo en el podio</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> con un increíble segundo lugar.</w:t>
      </w:r>
      <w:r>
        <w:br/>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Segoe UI Emoji" w:hAnsi="Segoe UI Emoji" w:cs="Segoe UI Emoji"/>
        </w:rPr>
        <w:t>🥉</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
        </w:rPr>
        <w:t>Raúl Tomás cerró el top 3</w:t>
      </w:r>
      <w:r>
        <w:t>, firmando una gran actuación.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Segoe UI Emoji" w:hAnsi="Segoe UI Emoji" w:cs="Segoe UI Emoji"/>
        </w:rPr>
        <w:t>💥</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> En la salida hubo caos: Tony Moreno trompeó tras ser frenado por Guti, lo que obligó a </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>Tejeriño</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> a entrar a boxes.</w:t>
      </w:r>
      <w:r>
        <w:br/>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Segoe UI Emoji" w:hAnsi="Segoe UI Emoji" w:cs="Segoe UI Emoji"/>
        </w:rPr>
        <w:t>💥</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> Armesto terminó fuera por un accidente.</w:t>
      </w:r>
      <w:r>
        <w:br/>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Segoe UI Emoji" w:hAnsi="Segoe UI Emoji" w:cs="Segoe UI Emoji"/>
        </w:rPr>
        <w:t>⚔️</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> Tony también tuvo un encontronazo con Axel Villar en plena lucha por posición.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Segoe UI Emoji" w:hAnsi="Segoe UI Emoji" w:cs="Segoe UI Emoji"/>
        </w:rPr>
        <w:t>🔍</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> La carrera no tuvo lluvia, pero sí mucha acción.</w:t>
      </w:r>
      <w:r>
        <w:br/>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Segoe UI Emoji" w:hAnsi="Segoe UI Emoji" w:cs="Segoe UI Emoji"/>
        </w:rPr>
        <w:t>📊</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> Los puntos fueron para: Cristian, Izan, </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>Shurde</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve">, Félix, </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>Tejeriño</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t>, Sánchez, Piñero y Guti.</w:t>
      </w:r>
      <w:r>
        <w:br/>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Segoe UI Emoji" w:hAnsi="Segoe UI Emoji" w:cs="Segoe UI Emoji"/>
        </w:rPr>
        <w:t>😬</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> Tony terminó 11°, </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
        </w:rPr>
        <w:t>a un paso de los puntos</w:t>
      </w:r>
      <w:r>
        <w:t>, después de haber clasificado 10°.</w:t>
      </w:r>
      <w:r>
        <w:br/>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Segoe UI Emoji" w:hAnsi="Segoe UI Emoji" w:cs="Segoe UI Emoji"/>
        </w:rPr>
        <w:t>❌</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> McLaren se fue con las manos vac</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/>
        </w:rPr>
        <w:t>í</w:t>
      </w:r>
      <w:r>
        <w:t>as.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Segoe UI Emoji" w:hAnsi="Segoe UI Emoji" w:cs="Segoe UI Emoji"/>
        </w:rPr>
        <w:t>📈</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> En el campeonato de pilotos:</w:t>
      </w:r>
      <w:r>
        <w:br/>
        <w:t>1️</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Segoe UI Symbol" w:hAnsi="Segoe UI Symbol" w:cs="Segoe UI Symbol"/>
        </w:rPr>
        <w:t>⃣</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> Diego lidera.</w:t>
      </w:r>
      <w:r>
        <w:br/>
        <w:t>2️</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Segoe UI Symbol" w:hAnsi="Segoe UI Symbol" w:cs="Segoe UI Symbol"/>
        </w:rPr>
        <w:t>⃣</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> Hugo pisa fuerte en su debut.</w:t>
      </w:r>
      <w:r>
        <w:br/>
        <w:t>3️</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Segoe UI Symbol" w:hAnsi="Segoe UI Symbol" w:cs="Segoe UI Symbol"/>
        </w:rPr>
        <w:t>⃣</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> Raúl Martín sigue sumando.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Segoe UI Emoji" w:hAnsi="Segoe UI Emoji" w:cs="Segoe UI Emoji"/>
        </w:rPr>
        <w:t>👑</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> En constructores:</w:t>
      </w:r>
      <w:r>
        <w:br/>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Segoe UI Emoji" w:hAnsi="Segoe UI Emoji" w:cs="Segoe UI Emoji"/>
        </w:rPr>
        <w:t>🥇</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>Kicker</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> al frente, seguido por Williams y Red Bull.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
        </w:rPr>
        <w:t>Silverstone</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
        </w:rPr>
        <w:t xml:space="preserve"> no decepciona. Diego manda, Hugo sorprende y el campeonato se aprieta.</w:t>
      </w:r>
    </w:p>
'@

$xmlPackage = @"
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
$newParagraphsXml
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@

# Find the last content paragraph (ends with the closing fire emoji),
# i.e. the paragraph immediately before the trailing empty paragraph.
$lastParaIndex = $d.Paragraphs.Count - 1
$target = $d.Paragraphs.Item($lastParaIndex).Range
$target.Collapse(0)
$target.InsertXML($xmlPackage)
